$wb = $excel.ActiveWorkbook

# --- Sheet "Runs": insert 4 new rows (Run112-Run115) before existing row 6 (Run116) ---
$runsWs = $wb.Worksheets.Item("Runs")

# Insert four blank rows at position 6; this pushes the existing row 6 (Run116) down to row 10
$runsWs.Rows.Item(6).Insert()
$runsWs.Rows.Item(6).Insert()
$runsWs.Rows.Item(6).Insert()
$runsWs.Rows.Item(6).Insert()

$runsNewRows = @(
    @{ Row = 6; Name = "Run112"; Values = @("0.357 (0.417)","0.561 (0.199)","0.282 (0.450)","0.304 (0.428)","0.720 (0.213)","0.210 (0.408)","0.456 (0.462)","0.593 (0.169)","0.429 (0.496)","0.381 (0.438)","0.594 (0.201)","0.321 (0.467)") },
    @{ Row = 7; Name = "Run113"; Values = @("0.232 (0.368)","0.581 (0.181)","0.143 (0.351)","0.305 (0.390)","0.574 (0.236)","0.180 (0.385)","0.361 (0.410)","0.568 (0.195)","0.272 (0.445)","0.292 (0.392)","0.574 (0.199)","0.193 (0.395)") },
    @{ Row = 8; Name = "Run114"; Values = @("0.456 (0.466)","0.648 (0.117)","0.423 (0.495)","0.345 (0.388)","0.601 (0.237)","0.165 (0.373)","0.588 (0.443)","0.775 (0.135)","0.511 (0.500)","0.499 (0.455)","0.703 (0.177)","0.432 (0.496)") },
    @{ Row = 9; Name = "Run115"; Values = @("0.515 (0.466)","0.620 (0.249)","0.490 (0.500)","0.558 (0.460)","0.769 (0.116)","0.500 (0.501)","0.672 (0.441)","0.529 (0.171)","0.692 (0.462)","0.580 (0.461)","0.630 (0.222)","0.570 (0.495)") }
)

foreach ($entry in $runsNewRows) {
    $r = $entry.Row
    $runsWs.Cells.Item($r, 1).Value = 0
    $runsWs.Cells.Item($r, 2).Value = $entry.Name
    $col = 3
    foreach ($val in $entry.Values) {
        $runsWs.Cells.Item($r, $col).Value = $val
        $col = $col + 1
    }
}

# --- Sheet "Iterations": add new row 2 (Iteration31) ---
$iterWs = $wb.Worksheets.Item("Iterations")

$iterWs.Cells.Item(2, 1).Value = 0
$iterWs.Cells.Item(2, 2).Value = "Iteration31"
$iterValues = @("0.39 (0.11)","0.60 (0.03)","0.33 (0.13)","0.38 (0.11)","0.67 (0.08)","0.26 (0.14)","0.52 (0.12)","0.62 (0.09)","0.48 (0.15)","0.44 (0.11)","0.63 (0.05)","0.38 (0.14)")
$col = 3
foreach ($val in $iterValues) {
    $iterWs.Cells.Item(2, $col).Value = $val
    $col = $col + 1
}
